# Add three new market sheets (Netherlands, Austria, Denmark) by copying the
# existing "Germany" template sheet to the end of the workbook, then editing
# the market name / SKU values (and, where needed, removing rows for
# products that are not sold in that market).

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("Germany")

# ---------------------------------------------------------------------------
# Netherlands
# ---------------------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$nl = $wb.Worksheets.Item($wb.Worksheets.Count)
$nl.Name = "Netherlands"
$nl.Range("B4").Value = "NGC-3144/T2188"
$nl.Range("B2").Value = "Netherlands Market"
$nl.Range("B4").Select()

# ---------------------------------------------------------------------------
# Austria (P32AR / P32DR rows are not applicable, so remove them)
# ---------------------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$at = $wb.Worksheets.Item($wb.Worksheets.Count)
$at.Name = "Austria"
$at.Range("B4").Value = "NGC-3817/T2295"
$at.Range("B2").Value = "Austria Market"
$at.Rows("16:17").Delete()
$at.Range("B4").Select()

# ---------------------------------------------------------------------------
# Denmark (custom "MZXSDR240" product replaces P32AR; P32DR row removed)
# ---------------------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$dk = $wb.Worksheets.Item($wb.Worksheets.Count)
$dk.Name = "Denmark"
$dk.Range("A16").Value = "MZXSDR240"
$dk.Range("B4").Value = "NGC-2913/T2783"
$dk.Range("B2").Value = "Denmark Market"
$dk.Rows(17).Delete()
$dk.Range("B4").Select()

# Austria ends up being the active/selected tab once everything is done.
$at.Activate()
